$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in homework grades (value 5) for the affected students
$ws.Range("C8:E8").Value = 5
$ws.Range("C20:F20").Value = 5
$ws.Range("C30:F30").Value = 5

# Update the frozen pane / scroll position and active selection
$ws.Range("F21").Select()
$excel.ActiveWindow.ScrollRow = 7
